$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<these>"
$ws.Range("C2").Value = 52

# Row 3
$ws.Range("C3").Value = 55

# Row 4
$ws.Range("B4").Value = "<if>"
$ws.Range("C4").Value = 56

# Row 5
$ws.Range("C5").Value = 53

# Row 6
$ws.Range("B6").Value = "<perte>"
$ws.Range("C6").Value = 51

# Row 7
$ws.Range("C7").Value = 54

# Row 8
$ws.Range("B8").Value = "<foxtrot>"
$ws.Range("C8").Value = 52

# Row 10
$ws.Range("C10").Value = 52

# Row 11
$ws.Range("B11").Value = "<we>"
$ws.Range("C11").Value = 52

# Row 12
$ws.Range("B12").Value = "<in>"
$ws.Range("C12").Value = 56

# Row 13
$ws.Range("B13").Value = "<oi>"
$ws.Range("C13").Value = 55

# Row 15
$ws.Range("B15").Value = "<more>"
$ws.Range("C15").Value = 60

# Row 16
$ws.Range("C16").Value = 30
